# Apply updated Price (D) and Volume(1h) (E) values from the latest cryptos refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.361.30"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "'3.457.08"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'576.30"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'160.81"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'3.458.21"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("E9").Value = "  +10.40%  "
$ws.Range("D10").Value = "'7.34"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "'4.055.62"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").Value = "'0.0000196"
$ws.Range("E15").Value = "  +6.57%  "
$ws.Range("D16").Value = "'29.06"
$ws.Range("E16").Value = "  +7.71%  "
$ws.Range("D17").Value = "'64.404.01"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "'3.449.72"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "'6.41"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "'14.53"
$ws.Range("E20").Value = "  +4.10%  "
$ws.Range("D21").Value = "'388.22"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'8.26"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "'0.549"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").Value = "'73.47"
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'0.0000125"
$ws.Range("E26").Value = "  +21.15%  "
$ws.Range("D27").Value = "'9.56"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'6.19"
$ws.Range("E30").Value = "  +10.78%  "
$ws.Range("E31").Value = "  +9.12%  "
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "'6.60"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'7.11"
$ws.Range("E36").Value = "  +5.96%  "
$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").Value = "'160.56"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("D40").Value = "'1.89"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'2.917.51"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "'4.53"
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("D45").Value = "'42.66"
$ws.Range("E45").Value = "  +3.85%  "
$ws.Range("D46").Value = "'0.774"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").Value = "'23.84"
$ws.Range("E47").Value = "  +7.77%  "
$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("D49").Value = "'2.24"
$ws.Range("E49").Value = "  +17.94%  "
$ws.Range("D50").Value = "'0.109"
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("D51").Value = "'0.873"
$ws.Range("E51").Value = "  +7.24%  "
